$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038736733008335
$ws.Range("D2").Value = 1.041381105907747
$ws.Range("E2").Value = 1.046682938841214
$ws.Range("F2").Value = 1.055942536798222
$ws.Range("I2").Value = 1.038222412079418
$ws.Range("J2").Value = 1.043832550664196
$ws.Range("K2").Value = 1.044160782336763
$ws.Range("L2").Value = 1.049447689178729
$ws.Range("M2").Value = 1.058681608422874
$ws.Range("N2").Value = 1.045314912615406

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039780827347481
$ws.Range("D3").Value = 1.042339731357351
$ws.Range("E3").Value = 1.047616995513887
$ws.Range("F3").Value = 1.056978071737368
$ws.Range("I3").Value = 1.038423802944622
$ws.Range("J3").Value = 1.044521087292708
$ws.Range("K3").Value = 1.04492963519491
$ws.Range("L3").Value = 1.050193119988467
$ws.Range("M3").Value = 1.059530118778651
$ws.Range("N3").Value = 1.046004427044906

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.040456342944165
$ws.Range("D4").Value = 1.042960212949667
$ws.Range("E4").Value = 1.048221675935279
$ws.Range("F4").Value = 1.057648404601008
$ws.Range("I4").Value = 1.038552131797294
$ws.Range("J4").Value = 1.044965956123034
$ws.Range("K4").Value = 1.045426702481787
$ws.Range("L4").Value = 1.050675121161455
$ws.Range("M4").Value = 1.060078827881878
$ws.Range("N4").Value = 1.04644992764

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.040740309766123
$ws.Range("D5").Value = 1.043221107631961
$ws.Range("E5").Value = 1.048475950844117
$ws.Range("F5").Value = 1.057930277276346
$ws.Range("I5").Value = 1.038605605736411
$ws.Range("J5").Value = 1.045152820496847
$ws.Range("K5").Value = 1.045635565746029
$ws.Range("L5").Value = 1.05087767204987
$ws.Range("M5").Value = 1.060309424792119
$ws.Range("N5").Value = 1.046637057382657

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.040787987875928
$ws.Range("D6").Value = 1.043264915575461
$ws.Range("E6").Value = 1.048518648647893
$ws.Range("F6").Value = 1.057977608720539
$ws.Range("I6").Value = 1.038614556349702
$ws.Range("J6").Value = 1.04518418654285
$ws.Range("K6").Value = 1.045670628706603
$ws.Range("L6").Value = 1.050911676390406
$ws.Range("M6").Value = 1.060348138299584
$ws.Range("N6").Value = 1.046668467972043

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.040460137401077
$ws.Range("D7").Value = 1.042963698864067
$ws.Range("E7").Value = 1.048225073305486
$ws.Range("F7").Value = 1.057652170743748
$ws.Range("I7").Value = 1.038552848187651
$ws.Range("J7").Value = 1.044968453637158
$ws.Range("K7").Value = 1.045429493731438
$ws.Range("L7").Value = 1.050677827980686
$ws.Range("M7").Value = 1.060081909445593
$ws.Range("N7").Value = 1.046452428700879

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.039089606817329
$ws.Range("D8").Value = 1.041705038881454
$ws.Range("E8").Value = 1.046998548672813
$ws.Range("F8").Value = 1.05629244395429
$ws.Range("I8").Value = 1.038290883721478
$ws.Range("J8").Value = 1.044065381357635
$ws.Range("K8").Value = 1.044420708942275
$ws.Range("L8").Value = 1.049699681650605
$ws.Range("M8").Value = 1.058968435308996
$ws.Range("N8").Value = 1.045548073955136

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036673918498246
$ws.Range("D9").Value = 1.039488571653383
$ws.Range("E9").Value = 1.044839449429905
$ws.Range("F9").Value = 1.053898540138887
$ws.Range("I9").Value = 1.03781408783029
$ws.Range("J9").Value = 1.042469010767774
$ws.Range("K9").Value = 1.042639811101048
$ws.Range("L9").Value = 1.047973458873449
$ws.Range("M9").Value = 1.057003815152195
$ws.Range("N9").Value = 1.043949436335931

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035063025895806
$ws.Range("D10").Value = 1.038011929723903
$ws.Range("E10").Value = 1.043401554926538
$ws.Range("F10").Value = 1.052304047707812
$ws.Range("I10").Value = 1.037486039725386
$ws.Range("J10").Value = 1.041401393099898
$ws.Range("K10").Value = 1.041450355278182
$ws.Range("L10").Value = 1.046820915393809
$ws.Range("M10").Value = 1.055692387655047
$ws.Range("N10").Value = 1.042880302528509

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034365386603556
$ws.Range("D11").Value = 1.037372769989
$ws.Range("E11").Value = 1.042779293179493
$ws.Range("F11").Value = 1.051613961719818
$ws.Range("I11").Value = 1.037341580035911
$ws.Range("J11").Value = 1.040938306668777
$ws.Range("K11").Value = 1.040934791906774
$ws.Range("L11").Value = 1.046321445377326
$ws.Range("M11").Value = 1.055124130743538
$ws.Range("N11").Value = 1.042416558461539

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034106235145569
$ws.Range("D12").Value = 1.03713539332951
$ws.Range("E12").Value = 1.042548211303983
$ws.Range("F12").Value = 1.051357684354243
$ws.Range("I12").Value = 1.037287559070557
$ws.Range("J12").Value = 1.040766175584577
$ws.Range("K12").Value = 1.040743210618652
$ws.Range("L12").Value = 1.046135858440937
$ws.Range("M12").Value = 1.054912994951077
$ws.Range("N12").Value = 1.042244182931456

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034161824794684
$ws.Range("D13").Value = 1.037186309807818
$ws.Range("E13").Value = 1.042597776708622
$ws.Range("F13").Value = 1.051412654393856
$ws.Range("I13").Value = 1.037299163137158
$ws.Range("J13").Value = 1.040803103721537
$ws.Range("K13").Value = 1.040784308986163
$ws.Range("L13").Value = 1.046175670249699
$ws.Range("M13").Value = 1.054958287000835
$ws.Range("N13").Value = 1.042281163510606

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034343965415115
$ws.Range("D14").Value = 1.037353147629943
$ws.Range("E14").Value = 1.042760190779185
$ws.Range("F14").Value = 1.051592776704965
$ws.Range("I14").Value = 1.037337122033353
$ws.Range("J14").Value = 1.04092408071004
$ws.Range("K14").Value = 1.040918957329739
$ws.Range("L14").Value = 1.046306105957236
$ws.Range("M14").Value = 1.055106679413087
$ws.Range("N14").Value = 1.042402312300309

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.034456185979218
$ws.Range("D15").Value = 1.037455946660448
$ws.Range("E15").Value = 1.042860266593548
$ws.Range("F15").Value = 1.051703762810342
$ws.Range("I15").Value = 1.037360461774625
$ws.Range("J15").Value = 1.040998602707305
$ws.Range("K15").Value = 1.041001908275533
$ws.Range("L15").Value = 1.046386463554519
$ws.Range("M15").Value = 1.055198100949183
$ws.Range("N15").Value = 1.042476940127358

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.035109323528358
$ws.Range("D16").Value = 1.038054353645479
$ws.Range("E16").Value = 1.043442859921072
$ws.Range("F16").Value = 1.052349853622001
$ws.Range("I16").Value = 1.037495576230139
$ws.Range("J16").Value = 1.041432109760754
$ws.Range("K16").Value = 1.041484560558364
$ws.Range("L16").Value = 1.046854054932292
$ws.Range("M16").Value = 1.055730092585014
$ws.Range("N16").Value = 1.042911062810547

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035518989209341
$ws.Range("D17").Value = 1.038429781907134
$ws.Range("E17").Value = 1.04380840089985
$ws.Range("F17").Value = 1.052755220460972
$ws.Range("I17").Value = 1.037579684207734
$ws.Range("J17").Value = 1.041703822968636
$ws.Range("K17").Value = 1.041787176200448
$ws.Range("L17").Value = 1.047147252753371
$ws.Range("M17").Value = 1.056063689925798
$ws.Range("N17").Value = 1.043183161882349

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035757929468834
$ws.Range("D18").Value = 1.038648785609608
$ws.Range("E18").Value = 1.044021649099997
$ws.Range("F18").Value = 1.052991696844933
$ws.Range("I18").Value = 1.037628510159797
$ws.Range("J18").Value = 1.041862231396216
$ws.Range("K18").Value = 1.041963636488
$ws.Range("L18").Value = 1.047318230403285
$ws.Range("M18").Value = 1.056258232985835
$ws.Range("N18").Value = 1.043341795268067

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035839400099874
$ws.Range("D19").Value = 1.03872346406283
$ws.Range("E19").Value = 1.044094367037952
$ws.Range("F19").Value = 1.053072334730534
$ws.Range("I19").Value = 1.037645119050886
$ws.Range("J19").Value = 1.041916231461573
$ws.Range("K19").Value = 1.042023796348545
$ws.Range("L19").Value = 1.047376522626269
$ws.Range("M19").Value = 1.056324560582901
$ws.Range("N19").Value = 1.043395872019711

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035475037078116
$ws.Range("D20").Value = 1.038389499652519
$ws.Range("E20").Value = 1.043769178258546
$ws.Range("F20").Value = 1.052711725044508
$ws.Range("I20").Value = 1.037570684296122
$ws.Range("J20").Value = 1.041674678706003
$ws.Range("K20").Value = 1.04175471360683
$ws.Range("L20").Value = 1.047115799504359
$ws.Range("M20").Value = 1.056027902081524
$ws.Range("N20").Value = 1.043153976231521

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034290330018414
$ws.Range("D21").Value = 1.037304017061306
$ws.Range("E21").Value = 1.042712362394344
$ws.Range("F21").Value = 1.051539733755745
$ws.Range("I21").Value = 1.037325954080958
$ws.Range("J21").Value = 1.040888459318019
$ws.Range("K21").Value = 1.04087930890284
$ws.Range("L21").Value = 1.046267697589211
$ws.Range("M21").Value = 1.05506298319714
$ws.Range("N21").Value = 1.042366640321826

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033545358216279
$ws.Range("D22").Value = 1.03662173789674
$ws.Range("E22").Value = 1.042048212006838
$ws.Range("F22").Value = 1.050803153222468
$ws.Range("I22").Value = 1.037169987009761
$ws.Range("J22").Value = 1.040393436826483
$ws.Range("K22").Value = 1.04032845547381
$ws.Range("L22").Value = 1.045734106350039
$ws.Range("M22").Value = 1.054455954102093
$ws.Range("N22").Value = 1.041870914841571

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03394029139564
$ws.Range("D23").Value = 1.036983407296744
$ws.Range("E23").Value = 1.042400261033502
$ws.Range("F23").Value = 1.051193600273864
$ws.Range("I23").Value = 1.037252866641107
$ws.Range("J23").Value = 1.040655923453987
$ws.Range("K23").Value = 1.040620516053918
$ws.Range("L23").Value = 1.046017006822825
$ws.Range("M23").Value = 1.054777784455464
$ws.Range("N23").Value = 1.042133774230194

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03549489718748
$ws.Range("D24").Value = 1.038407701399062
$ws.Range("E24").Value = 1.043786901175676
$ws.Range("F24").Value = 1.052731378649621
$ws.Range("I24").Value = 1.037574751688086
$ws.Range("J24").Value = 1.0416878479825
$ws.Range("K24").Value = 1.04176938221086
$ws.Range("L24").Value = 1.047130011996211
$ws.Range("M24").Value = 1.056044073187125
$ws.Range("N24").Value = 1.043167164209902

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037298507928588
$ws.Range("D25").Value = 1.040061406697026
$ws.Range("E25").Value = 1.045397365242152
$ws.Range("F25").Value = 1.054517169103741
$ws.Range("I25").Value = 1.037939147665924
$ws.Range("J25").Value = 1.04288230630867
$ws.Range("K25").Value = 1.043100603539905
$ws.Range("L25").Value = 1.048420035085152
$ws.Range("M25").Value = 1.057512014573341
$ws.Range("N25").Value = 1.044363318803901

